$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 7500
$ws.Range("J46").Value = 6666.6665
$ws.Range("L46").Value = 19999.9995
$ws.Range("N46").Value = -20237.9995
$ws.Range("H60").Value = 7500
$ws.Range("J60").Value = 6666.6665
$ws.Range("L60").Value = 19999.9995
$ws.Range("N60").Value = -20967.9995
$ws.Range("H133").Value = 42810.555
$ws.Range("J133").Value = 42810.555
$ws.Range("L133").Value = 42810.555
$ws.Range("N133").Value = -52930.555
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22746.7
$ws.Range("I32").Value = 4851.9805
$ws.Range("J32").Value = 70779.89999999999
$ws.Range("K32").Value = 4851.9805
$ws.Range("L32").Value = 70779.89999999999
$ws.Range("M32").Value = -4564.9805
$ws.Range("N32").Value = -71353.89999999999
$ws.Range("H45").Value = 995
$ws.Range("I45").Value = 995
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 995
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -618
$ws.Range("H132").Value = 4214.56
$ws.Range("I132").Value = 4235.3335
$ws.Range("J132").Value = 4195.385
$ws.Range("K132").Value = 12706.0005
$ws.Range("L132").Value = 12586.155
$ws.Range("M132").Value = -10176.0005
$ws.Range("N132").Value = -17646.155
$ws.Range("H133").Value = 48250
$ws.Range("J133").Value = 48250
$ws.Range("L133").Value = 48250
$ws.Range("N133").Value = -53310
$ws.Range("H139").Value = 53571.668
$ws.Range("J139").Value = 53571.668
$ws.Range("L139").Value = 53571.668
$ws.Range("N139").Value = -63851.668
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 13410
$ws.Range("I5").Value = 10174
$ws.Range("J5").Value = 21500
$ws.Range("K5").Value = 10174
$ws.Range("L5").Value = 21500
$ws.Range("M5").Value = -10061
$ws.Range("N5").Value = -21726
$ws.Range("H59").Value = 55450
$ws.Range("J59").Value = 55450
$ws.Range("L59").Value = 55450
$ws.Range("N59").Value = -57144
$ws.Range("H86").Value = 1521.6923
$ws.Range("I86").Value = 1535.6
$ws.Range("J86").Value = 1513
$ws.Range("K86").Value = 1535.6
$ws.Range("L86").Value = 1513
$ws.Range("M86").Value = -412.5999999999999
$ws.Range("N86").Value = -3759
$ws.Range("H89").Value = 1521.6923
$ws.Range("I89").Value = 1535.6
$ws.Range("J89").Value = 1513
$ws.Range("K89").Value = 7678
$ws.Range("L89").Value = 7565
$ws.Range("M89").Value = -2062
$ws.Range("N89").Value = -18797
$ws.Range("H105").Value = 235562.33
$ws.Range("I105").Value = 2802.3794
$ws.Range("K105").Value = 2802.3794
$ws.Range("M105").Value = -1055.3794
$ws.Range("H133").Value = 41390
$ws.Range("J133").Value = 41390
$ws.Range("L133").Value = 41390
$ws.Range("N133").Value = -51510
$ws.Range("H134").Value = 2148.94
$ws.Range("I134").Value = 1556.3334
$ws.Range("J134").Value = 5260.125
$ws.Range("K134").Value = 4669.0002
$ws.Range("L134").Value = 15780.375
$ws.Range("M134").Value = -2134.0002
$ws.Range("N134").Value = -20850.375
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 50000
$ws.Range("I2").Value = 50000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 50000
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -49887
$ws.Range("H31").Value = 4651.491
$ws.Range("I31").Value = 1395.8857
$ws.Range("J31").Value = 9830.862999999999
$ws.Range("K31").Value = 1395.8857
$ws.Range("L31").Value = 9830.862999999999
$ws.Range("M31").Value = -1100.8857
$ws.Range("N31").Value = -10420.863
$ws.Range("H34").Value = 4651.491
$ws.Range("I34").Value = 1395.8857
$ws.Range("J34").Value = 9830.862999999999
$ws.Range("K34").Value = 1395.8857
$ws.Range("L34").Value = 9830.862999999999
$ws.Range("M34").Value = -1193.8857
$ws.Range("N34").Value = -10234.863
$ws.Range("H134").Value = 5035.364
$ws.Range("I134").Value = 1934
$ws.Range("J134").Value = 7619.8335
$ws.Range("K134").Value = 5802
$ws.Range("L134").Value = 22859.5005
$ws.Range("M134").Value = -3267
$ws.Range("N134").Value = -27929.5005
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1741.2195
$ws.Range("I5").Value = 1477.8572
$ws.Range("J5").Value = 1877.7778
$ws.Range("K5").Value = 4433.571599999999
$ws.Range("L5").Value = 5633.3334
$ws.Range("M5").Value = -4321.571599999999
$ws.Range("N5").Value = -5857.3334
$ws.Range("H61").Value = 156.8
$ws.Range("I61").Value = 156.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 470.4
$ws.Range("L61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -255.4
$ws.Range("H131").Value = 1302.6154
$ws.Range("I131").Value = 400
$ws.Range("J131").Value = 1405.7715
$ws.Range("K131").Value = 1200
$ws.Range("L131").Value = 4217.3145
$ws.Range("M131").Value = 3840
$ws.Range("N131").Value = -14297.3145
$ws.Range("H135").Value = 1741.2195
$ws.Range("I135").Value = 1477.8572
$ws.Range("J135").Value = 1877.7778
$ws.Range("K135").Value = 13300.7148
$ws.Range("L135").Value = 16900.0002
$ws.Range("M135").Value = -10765.7148
$ws.Range("N135").Value = -21970.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 50000
$ws.Range("J4").Value = 50000
$ws.Range("L4").Value = 50000
$ws.Range("N4").Value = -50224
$ws.Range("H70").Value = 5288.5137
$ws.Range("I70").Value = 5378.9653
$ws.Range("K70").Value = 5378.9653
$ws.Range("M70").Value = -5108.9653
$ws.Range("H73").Value = 5288.5137
$ws.Range("I73").Value = 5378.9653
$ws.Range("K73").Value = 5378.9653
$ws.Range("M73").Value = -4442.9653
$ws.Range("H132").Value = 3316.9092
$ws.Range("I132").Value = 2773.4211
$ws.Range("J132").Value = 4054.5
$ws.Range("K132").Value = 8320.263300000001
$ws.Range("L132").Value = 12163.5
$ws.Range("M132").Value = -5790.263300000001
$ws.Range("N132").Value = -17223.5
$ws.Range("H137").Value = 48560
$ws.Range("J137").Value = 48560
$ws.Range("L137").Value = 48560
$ws.Range("N137").Value = -58760
$ws.Range("H138").Value = 77571.5
$ws.Range("J138").Value = 77571.5
$ws.Range("L138").Value = 77571.5
$ws.Range("N138").Value = -87851.5
$ws.Range("H139").Value = 42828
$ws.Range("J139").Value = 42828
$ws.Range("L139").Value = 42828
$ws.Range("N139").Value = -53108
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 747.375
$ws.Range("I46").Value = 613.25
$ws.Range("K46").Value = 613.25
$ws.Range("M46").Value = -425.25
$ws.Range("H132").Value = 3400.3845
$ws.Range("I132").Value = 1960.7142
$ws.Range("J132").Value = 7065
$ws.Range("K132").Value = 5882.142599999999
$ws.Range("L132").Value = 21195
$ws.Range("M132").Value = -3352.142599999999
$ws.Range("N132").Value = -26255
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1810.5278
$ws.Range("I132").Value = 1276.8148
$ws.Range("J132").Value = 3411.6667
$ws.Range("K132").Value = 3830.4444
$ws.Range("L132").Value = 10235.0001
$ws.Range("M132").Value = -1300.4444
$ws.Range("N132").Value = -15295.0001
